$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 169.17
$ws.Range("I15").Value = 169.17
$ws.Range("K15").Value = 507.51
$ws.Range("M15").Value = -338.51
$ws.Range("H17").Value = 1460.5454
$ws.Range("J17").Value = 1509.0476
$ws.Range("L17").Value = 4527.142800000001
$ws.Range("N17").Value = -4863.142800000001
$ws.Range("H112").Value = 10418229
$ws.Range("I112").Value = 500000350
$ws.Range("J112").Value = 1588.1277
$ws.Range("K112").Value = 1500001050
$ws.Range("L112").Value = 4764.3831
$ws.Range("M112").Value = -1499999942
$ws.Range("N112").Value = -6980.3831
$ws.Range("H138").Value = 2332.53
$ws.Range("I138").Value = 958.3103599999999
$ws.Range("J138").Value = 2893.831
$ws.Range("K138").Value = 2874.93108
$ws.Range("L138").Value = 8681.493
$ws.Range("M138").Value = 2265.06892
$ws.Range("N138").Value = -18961.493
$ws.Range("H141").Value = 6499.3184
$ws.Range("I141").Value = 6623.095
$ws.Range("J141").Value = 3900
$ws.Range("K141").Value = 19869.285
$ws.Range("L141").Value = 11700
$ws.Range("M141").Value = -14689.285
$ws.Range("N141").Value = -22060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 9544.286
$ws.Range("I3").Value = 5852.5
$ws.Range("J3").Value = 11021
$ws.Range("K3").Value = 5852.5
$ws.Range("L3").Value = 11021
$ws.Range("M3").Value = -5737.5
$ws.Range("N3").Value = -11251
$ws.Range("H32").Value = 7284.3286
$ws.Range("I32").Value = 4411.8604
$ws.Range("J32").Value = 11859
$ws.Range("K32").Value = 4411.8604
$ws.Range("L32").Value = 11859
$ws.Range("M32").Value = -4124.8604
$ws.Range("N32").Value = -12433
$ws.Range("H61").Value = 2150.1035
$ws.Range("I61").Value = 1715.5264
$ws.Range("J61").Value = 2975.8
$ws.Range("K61").Value = 1715.5264
$ws.Range("L61").Value = 2975.8
$ws.Range("M61").Value = -1503.5264
$ws.Range("N61").Value = -3399.8
$ws.Range("H74").Value = 8889.643
$ws.Range("I74").Value = 9959.5
$ws.Range("J74").Value = 6215
$ws.Range("K74").Value = 9959.5
$ws.Range("L74").Value = 6215
$ws.Range("M74").Value = -9085.5
$ws.Range("N74").Value = -7963
$ws.Range("H77").Value = 8889.643
$ws.Range("I77").Value = 9959.5
$ws.Range("J77").Value = 6215
$ws.Range("K77").Value = 49797.5
$ws.Range("L77").Value = 31075
$ws.Range("M77").Value = -45429.5
$ws.Range("N77").Value = -39811
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680
$ws.Range("H133").Value = 37920
$ws.Range("J133").Value = 37920
$ws.Range("L133").Value = 37920
$ws.Range("N133").Value = -42980
$ws.Range("H136").Value = 2150.1035
$ws.Range("I136").Value = 1715.5264
$ws.Range("J136").Value = 2975.8
$ws.Range("K136").Value = 5146.5792
$ws.Range("L136").Value = 8927.400000000001
$ws.Range("M136").Value = -2596.5792
$ws.Range("N136").Value = -14027.4
$ws.Range("H137").Value = 47690
$ws.Range("J137").Value = 47690
$ws.Range("L137").Value = 47690
$ws.Range("N137").Value = -57890

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1777.7646
$ws.Range("I20").Value = 1286.3077
$ws.Range("J20").Value = 3375
$ws.Range("K20").Value = 1286.3077
$ws.Range("L20").Value = 3375
$ws.Range("M20").Value = -1039.3077
$ws.Range("N20").Value = -3869
$ws.Range("H109").Value = 35000
$ws.Range("I109").Value = 35000
$ws.Range("K109").Value = 35000
$ws.Range("M109").Value = -33613
$ws.Range("H132").Value = 49902.11
$ws.Range("J132").Value = 49902.11
$ws.Range("L132").Value = 49902.11
$ws.Range("N132").Value = -60022.11
$ws.Range("H134").Value = 2582.7693
$ws.Range("I134").Value = 1521.7073
$ws.Range("J134").Value = 6537.636
$ws.Range("K134").Value = 4565.1219
$ws.Range("L134").Value = 19612.908
$ws.Range("M134").Value = -2030.1219
$ws.Range("N134").Value = -24682.908
$ws.Range("H137").Value = 39780
$ws.Range("J137").Value = 39780
$ws.Range("L137").Value = 39780
$ws.Range("N137").Value = -49980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1792.359
$ws.Range("I58").Value = 1636.3934
$ws.Range("J58").Value = 2352
$ws.Range("K58").Value = 1636.3934
$ws.Range("L58").Value = 2352
$ws.Range("M58").Value = -1433.3934
$ws.Range("N58").Value = -2758
$ws.Range("H136").Value = 1792.359
$ws.Range("I136").Value = 1636.3934
$ws.Range("J136").Value = 2352
$ws.Range("K136").Value = 4909.1802
$ws.Range("L136").Value = 7056
$ws.Range("M136").Value = -2359.1802
$ws.Range("N136").Value = -12156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 583.4706
$ws.Range("I113").Value = 563.6
$ws.Range("J113").Value = 638.6667
$ws.Range("K113").Value = 1690.8
$ws.Range("L113").Value = 1916.0001
$ws.Range("M113").Value = 479.1999999999998
$ws.Range("N113").Value = -6256.0001
$ws.Range("H131").Value = 5817108
$ws.Range("J131").Value = 755.7848
$ws.Range("L131").Value = 2267.3544
$ws.Range("N131").Value = -12347.3544

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 30838.572
$ws.Range("J46").Value = 30838.572
$ws.Range("L46").Value = 30838.572
$ws.Range("N46").Value = -31150.572
$ws.Range("H70").Value = 6277.2144
$ws.Range("I70").Value = 5605.8486
$ws.Range("K70").Value = 5605.8486
$ws.Range("M70").Value = -5335.8486
$ws.Range("H73").Value = 6277.2144
$ws.Range("I73").Value = 5605.8486
$ws.Range("K73").Value = 5605.8486
$ws.Range("M73").Value = -4669.8486
$ws.Range("H80").Value = 13891422
$ws.Range("J80").Value = 2760
$ws.Range("L80").Value = 2760
$ws.Range("N80").Value = -4756
$ws.Range("H83").Value = 13891422
$ws.Range("J83").Value = 2760
$ws.Range("L83").Value = 13800
$ws.Range("N83").Value = -23784
$ws.Range("H132").Value = 2945.7693
$ws.Range("I132").Value = 1805.421
$ws.Range("K132").Value = 5416.263
$ws.Range("M132").Value = -2886.263
$ws.Range("H137").Value = 69234.5
$ws.Range("J137").Value = 69234.5
$ws.Range("L137").Value = 69234.5
$ws.Range("N137").Value = -79434.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6433.222
$ws.Range("I7").Value = 2626
$ws.Range("J7").Value = 9479
$ws.Range("K7").Value = 2626
$ws.Range("L7").Value = 9479
$ws.Range("M7").Value = -2514
$ws.Range("N7").Value = -9703
$ws.Range("H100").Value = 2232.5557
$ws.Range("I100").Value = 1765.5
$ws.Range("J100").Value = 3166.6667
$ws.Range("K100").Value = 1765.5
$ws.Range("L100").Value = 3166.6667
$ws.Range("M100").Value = -1224.5
$ws.Range("N100").Value = -4248.6667
$ws.Range("H126").Value = 6433.222
$ws.Range("I126").Value = 2626
$ws.Range("J126").Value = 9479
$ws.Range("K126").Value = 7878
$ws.Range("L126").Value = 28437
$ws.Range("M126").Value = -5408
$ws.Range("N126").Value = -33377
$ws.Range("H132").Value = 3906.4126
$ws.Range("I132").Value = 1883.2222
$ws.Range("J132").Value = 8964.388999999999
$ws.Range("K132").Value = 5649.6666
$ws.Range("L132").Value = 26893.167
$ws.Range("M132").Value = -3119.6666
$ws.Range("N132").Value = -31953.167
$ws.Range("H136").Value = 3998.1538
$ws.Range("I136").Value = 2020.4117
$ws.Range("J136").Value = 7733.8887
$ws.Range("K136").Value = 6061.2351
$ws.Range("L136").Value = 23201.6661
$ws.Range("M136").Value = -3511.2351
$ws.Range("N136").Value = -28301.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 30999.75
$ws.Range("I49").Value = 10000
$ws.Range("J49").Value = 37999.668
$ws.Range("K49").Value = 10000
$ws.Range("L49").Value = 37999.668
$ws.Range("M49").Value = -9770
$ws.Range("N49").Value = -38459.668
$ws.Range("H122").Value = 3117.3428
$ws.Range("I122").Value = 2056.8262
$ws.Range("J122").Value = 5150
$ws.Range("K122").Value = 6170.4786
$ws.Range("L122").Value = 15450
$ws.Range("M122").Value = -3720.4786
$ws.Range("N122").Value = -20350
$ws.Range("H132").Value = 19611004
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 20836628
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 62509884
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -62514944
$ws.Range("H136").Value = 2180.2856
$ws.Range("I136").Value = 1385.9166
$ws.Range("J136").Value = 3239.4443
$ws.Range("K136").Value = 4157.7498
$ws.Range("L136").Value = 9718.332900000001
$ws.Range("M136").Value = -1607.7498
$ws.Range("N136").Value = -14818.3329
